# Updates the hourly crypto snapshot in Sheet1 (Coin / Link / Price / Volume(1h))
# with refreshed D (Price) and E (Volume 1h %) figures, matching the
# "Updated cryptos list ... with GitHub Actions" scrape refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "6.20") but must
# stay stored as text (matching the source price-string formatting, incl.
# any trailing zeros). Mark them as Text *before* writing so Excel does not
# silently reinterpret the string as a numeric value.
$textCells = @(
    "D5",
    "D6",
    "D11",
    "D14",
    "D20",
    "D22",
    "D23",
    "D25",
    "D32",
    "D34",
    "D36",
    "D37",
    "D45",
    "D46",
    "D51",
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Refreshed crypto price / 1h-volume-change figures.
$ws.Range("D2").Value = "64.795.21"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "3.160.09"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "573.19"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "151.12"
$ws.Range("E6").Value = "  +5.80%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.156.46"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("E9").Value = "  +5.07%  "
$ws.Range("E10").Value = "  +7.13%  "
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("E12").Value = "  +7.89%  "
$ws.Range("E13").Value = "  +13.20%  "
$ws.Range("D14").Value = "37.96"
$ws.Range("E14").Value = "  +8.86%  "
$ws.Range("D15").Value = "3.675.34"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "64.928.14"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("E17").Value = "  +8.02%  "
$ws.Range("D18").Value = "3.158.50"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D20").Value = "518.36"
$ws.Range("E20").Value = "  +8.22%  "
$ws.Range("E21").Value = "  +8.21%  "
$ws.Range("D22").Value = "0.737"
$ws.Range("E22").Value = "  +10.22%  "
$ws.Range("D23").Value = "15.18"
$ws.Range("E23").Value = "  +8.97%  "
$ws.Range("D25").Value = "85.28"
$ws.Range("E25").Value = "  +5.32%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +5.55%  "
$ws.Range("E28").Value = "  +10.58%  "
$ws.Range("E29").Value = "  +6.98%  "
$ws.Range("E30").Value = "  +6.91%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  +7.96%  "
$ws.Range("D34").Value = "6.12"
$ws.Range("E34").Value = "  +10.35%  "
$ws.Range("E35").Value = "  +7.07%  "
$ws.Range("D36").Value = "55.80"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "487.75"
$ws.Range("E37").Value = "  +9.07%  "
$ws.Range("E38").Value = "  +6.22%  "
$ws.Range("E39").Value = "  +4.66%  "
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").Value = "3.115.82"
$ws.Range("E41").Value = "  +5.14%  "
$ws.Range("E42").Value = "  +5.68%  "
$ws.Range("E43").Value = "  +6.31%  "
$ws.Range("E44").Value = "  +14.53%  "
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  +17.49%  "
$ws.Range("D46").Value = "29.21"
$ws.Range("E46").Value = "  +5.56%  "
$ws.Range("E47").Value = "  +14.18%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("E50").Value = "  +11.04%  "
$ws.Range("D51").Value = "118.84"
$ws.Range("E51").Value = "  -0.32%  "

